$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 111
$c = $ws.Cells.Item(111, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(111, 2).Value = 'v3'
$ws.Cells.Item(111, 3).Value = 'Day/'
$ws.Cells.Item(111, 4).Value = 2986
$ws.Cells.Item(111, 5).Value = 65
$ws.Cells.Item(111, 6).Value = 46
$ws.Cells.Item(111, 7).Value = 65
$ws.Cells.Item(111, 8).Value = 456
$ws.Cells.Item(111, 9).Value = 4
$ws.Cells.Item(111, 10).Value = 2990
$ws.Cells.Item(111, 11).Value = 'Nothing'
$ws.Cells.Item(111, 12).Value = 456
$ws.Cells.Item(111, 13).Value = 710
$c = $ws.Cells.Item(111, 14)
$c.NumberFormat = "@"
$c.Value = '456'
$c.ClearFormats()

# Row 112
$c = $ws.Cells.Item(112, 1)
$c.NumberFormat = "@"
$c.Value = '06-03-2018'
$c.ClearFormats()
$ws.Cells.Item(112, 2).Value = 'v1'
$ws.Cells.Item(112, 3).Value = 'Day/'
$ws.Cells.Item(112, 4).Value = 1624
$ws.Cells.Item(112, 5).Value = 3666
$ws.Cells.Item(112, 6).Value = 5
$ws.Cells.Item(112, 7).Value = 3666
$ws.Cells.Item(112, 8).Value = 87
$ws.Cells.Item(112, 9).Value = 7
$ws.Cells.Item(112, 10).Value = 18330
$ws.Cells.Item(112, 11).Value = 'Nothing'
$ws.Cells.Item(112, 12).Value = 74
$ws.Cells.Item(112, 13).Value = 17647
$c = $ws.Cells.Item(112, 14)
$c.NumberFormat = "@"
$c.Value = '7'
$c.ClearFormats()

# Row 113
$c = $ws.Cells.Item(113, 1)
$c.NumberFormat = "@"
$c.Value = '08-03-2018'
$c.ClearFormats()
$ws.Cells.Item(113, 2).Value = 'v3'
$ws.Cells.Item(113, 3).Value = '/Night'
$ws.Cells.Item(113, 4).Value = 65
$ws.Cells.Item(113, 5).Value = 76
$ws.Cells.Item(113, 6).Value = 456
$ws.Cells.Item(113, 7).Value = 76
$ws.Cells.Item(113, 8).Value = 456
$ws.Cells.Item(113, 9).Value = 456
$ws.Cells.Item(113, 10).Value = 34656
$ws.Cells.Item(113, 11).Value = 'Nothing'
$ws.Cells.Item(113, 12).Value = 46
$ws.Cells.Item(113, 13).Value = -173326
$c = $ws.Cells.Item(113, 14)
$c.NumberFormat = "@"
$c.Value = '465'
$c.ClearFormats()

# Row 114
$ws.Cells.Item(114, 1).Value = '13-03-2018'
$ws.Cells.Item(114, 2).Value = 'v2'
$ws.Cells.Item(114, 3).Value = 'Day/'
$ws.Cells.Item(114, 4).Value = 1650
$ws.Cells.Item(114, 5).Value = 4536
$ws.Cells.Item(114, 6).Value = 4563
$ws.Cells.Item(114, 7).Value = 4536
$ws.Cells.Item(114, 8).Value = 456
$ws.Cells.Item(114, 9).Value = 456
$ws.Cells.Item(114, 10).Value = 20697768
$ws.Cells.Item(114, 11).Value = 'Nothing'
$ws.Cells.Item(114, 12).Value = 465
$ws.Cells.Item(114, 13).Value = 20489367
$c = $ws.Cells.Item(114, 14)
$c.NumberFormat = "@"
$c.Value = '45'
$c.ClearFormats()

# Row 115
$c = $ws.Cells.Item(115, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(115, 2).Value = 'v1'
$ws.Cells.Item(115, 3).Value = '/Night'
$ws.Cells.Item(115, 4).Value = 3666
$ws.Cells.Item(115, 5).Value = 5
$ws.Cells.Item(115, 6).Value = 54
$ws.Cells.Item(115, 7).Value = 5
$ws.Cells.Item(115, 8).Value = 54
$ws.Cells.Item(115, 9).Value = 54
$ws.Cells.Item(115, 10).Value = 270
$ws.Cells.Item(115, 11).Value = 'Nothing'
$ws.Cells.Item(115, 12).Value = 45
$ws.Cells.Item(115, 13).Value = -2691
$c = $ws.Cells.Item(115, 14)
$c.NumberFormat = "@"
$c.Value = '45'
$c.ClearFormats()

# Row 116
$c = $ws.Cells.Item(116, 1)
$c.NumberFormat = "@"
$c.Value = '05-03-2018'
$c.ClearFormats()
$ws.Cells.Item(116, 2).Value = 'v1'
$ws.Cells.Item(116, 3).Value = '/Night'
$ws.Cells.Item(116, 4).Value = 5
$ws.Cells.Item(116, 5).Value = 5
$ws.Cells.Item(116, 6).Value = 456
$ws.Cells.Item(116, 7).Value = 5
$ws.Cells.Item(116, 8).Value = 45
$ws.Cells.Item(116, 9).Value = 45
$ws.Cells.Item(116, 10).Value = 2280
$ws.Cells.Item(116, 11).Value = 'Nothing'
$ws.Cells.Item(116, 12).Value = 45
$ws.Cells.Item(116, 13).Value = 210

# Row 117
$c = $ws.Cells.Item(117, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(117, 2).Value = 'v2'
$ws.Cells.Item(117, 3).Value = 'Day/'
$ws.Cells.Item(117, 4).Value = 4536
$ws.Cells.Item(117, 5).Value = 456
$ws.Cells.Item(117, 6).Value = 456
$ws.Cells.Item(117, 7).Value = 456
$ws.Cells.Item(117, 8).Value = 456468
$ws.Cells.Item(117, 9).Value = 468
$ws.Cells.Item(117, 10).Value = 207936
$ws.Cells.Item(117, 11).Value = 'Nothing'
$ws.Cells.Item(117, 12).Value = 486
$ws.Cells.Item(117, 13).Value = -213419574

# Row 118
$ws.Cells.Item(118, 1).Value = '23-03-2018'
$ws.Cells.Item(118, 2).Value = 'v1'
$ws.Cells.Item(118, 3).Value = 'Day/'
$ws.Cells.Item(118, 4).Value = 5
$ws.Cells.Item(118, 5).Value = 56
$ws.Cells.Item(118, 6).Value = 56
$ws.Cells.Item(118, 7).Value = 56
$ws.Cells.Item(118, 8).Value = 56
$ws.Cells.Item(118, 9).Value = 56
$ws.Cells.Item(118, 10).Value = 3136
$ws.Cells.Item(118, 11).Value = 'Nothing'
$ws.Cells.Item(118, 12).Value = 56
$ws.Cells.Item(118, 13).Value = -56

# Row 119
$ws.Cells.Item(119, 1).Value = '23-03-2018'
$ws.Cells.Item(119, 2).Value = 'v1'
$ws.Cells.Item(119, 3).Value = 'Day/Night'
$ws.Cells.Item(119, 4).Value = 5
$ws.Cells.Item(119, 5).Value = 46
$ws.Cells.Item(119, 6).Value = 4564
$ws.Cells.Item(119, 7).Value = 46
$ws.Cells.Item(119, 8).Value = 56
$ws.Cells.Item(119, 9).Value = 456
$ws.Cells.Item(119, 10).Value = 209944
$ws.Cells.Item(119, 11).Value = 'Nothing'
$ws.Cells.Item(119, 12).Value = 465
$ws.Cells.Item(119, 13).Value = 183943
$ws.Cells.Item(119, 14).Value = 'issued'

# Row 120
$ws.Cells.Item(120, 1).Value = '23-03-2018'
$ws.Cells.Item(120, 2).Value = 'v3'
$ws.Cells.Item(120, 3).Value = '/Night'
$ws.Cells.Item(120, 4).Value = 76
$ws.Cells.Item(120, 5).Value = 564
$ws.Cells.Item(120, 6).Value = 654
$ws.Cells.Item(120, 7).Value = 564
$ws.Cells.Item(120, 8).Value = 64
$ws.Cells.Item(120, 9).Value = 64
$ws.Cells.Item(120, 10).Value = 368856
$ws.Cells.Item(120, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = 364760
$ws.Cells.Item(120, 14).Value = 'issued'

# Row 121
$ws.Cells.Item(121, 1).Value = '16-03-2018'
$ws.Cells.Item(121, 2).Value = 'v1'
$ws.Cells.Item(121, 3).Value = 'Day/'
$ws.Cells.Item(121, 4).Value = 46
$ws.Cells.Item(121, 5).Value = 35
$ws.Cells.Item(121, 6).Value = 65
$ws.Cells.Item(121, 7).Value = 35
$ws.Cells.Item(121, 8).Value = 651
$ws.Cells.Item(121, 9).Value = 561
$ws.Cells.Item(121, 10).Value = 2275
$ws.Cells.Item(121, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = -362936
$ws.Cells.Item(121, 14).Value = 'issued'

# Row 122
$c = $ws.Cells.Item(122, 1)
$c.NumberFormat = "@"
$c.Value = '01-03-2018'
$c.ClearFormats()
$ws.Cells.Item(122, 2).Value = 'v1'
$ws.Cells.Item(122, 3).Value = '/Night'
$ws.Cells.Item(122, 4).Value = 35
$ws.Cells.Item(122, 5).Value = 564
$ws.Cells.Item(122, 6).Value = 654
$ws.Cells.Item(122, 7).Value = 564
$ws.Cells.Item(122, 8).Value = 64
$ws.Cells.Item(122, 9).Value = 65465
$ws.Cells.Item(122, 10).Value = 368856
$ws.Cells.Item(122, 11).Value = 'Air filter                 250'
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -3820904
$ws.Cells.Item(122, 14).Value = 'issued'

# Row 123
$ws.Cells.Item(123, 1).Value = '15-03-2018'
$ws.Cells.Item(123, 2).Value = 'v1'
$ws.Cells.Item(123, 3).Value = 'Day/'
$ws.Cells.Item(123, 4).Value = 564
$ws.Cells.Item(123, 5).Value = 54
$ws.Cells.Item(123, 6).Value = 534
$ws.Cells.Item(123, 7).Value = 54
$ws.Cells.Item(123, 8).Value = 54
$ws.Cells.Item(123, 9).Value = 564
$ws.Cells.Item(123, 10).Value = 28836
$ws.Cells.Item(123, 11).Value = 'Track motor oil       1000'
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).Value = -1620
$ws.Cells.Item(123, 14).Value = 'issued'

# Row 124
$c = $ws.Cells.Item(124, 1)
$c.NumberFormat = "@"
$c.Value = '06-03-2018'
$c.ClearFormats()
$ws.Cells.Item(124, 2).Value = 'v1'
$ws.Cells.Item(124, 3).Value = '/Night'
$ws.Cells.Item(124, 4).Value = 54
$ws.Cells.Item(124, 5).Value = 84
$ws.Cells.Item(124, 6).Value = 98
$ws.Cells.Item(124, 7).Value = 84
$ws.Cells.Item(124, 8).Value = 898
$ws.Cells.Item(124, 9).Value = 684
$ws.Cells.Item(124, 10).Value = 8232
$ws.Cells.Item(124, 11).Value = 'Hydraulic strainer  250'
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 13).Value = -606000
$ws.Cells.Item(124, 14).Value = 'issued'

# Row 125
$ws.Cells.Item(125, 1).Value = '22-03-2018'
$ws.Cells.Item(125, 2).Value = 'v2'
$ws.Cells.Item(125, 3).Value = 'Day/'
$ws.Cells.Item(125, 4).Value = 4536
$ws.Cells.Item(125, 5).Value = 451
$ws.Cells.Item(125, 6).Value = 654
$ws.Cells.Item(125, 7).Value = 451
$ws.Cells.Item(125, 8).Value = 514
$ws.Cells.Item(125, 9).Value = 165
$ws.Cells.Item(125, 10).Value = 294954
$ws.Cells.Item(125, 11).Value = 'Nothing'
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = 210144
$ws.Cells.Item(125, 14).Value = 'issued'

# Row 126
$ws.Cells.Item(126, 1).Value = '22-03-2018'
$ws.Cells.Item(126, 2).Value = 'v2'
$ws.Cells.Item(126, 3).Value = 'Day/'
$ws.Cells.Item(126, 4).Value = 451
$ws.Cells.Item(126, 5).Value = 451
$ws.Cells.Item(126, 6).Value = 654
$ws.Cells.Item(126, 7).Value = 451
$ws.Cells.Item(126, 8).Value = 514
$ws.Cells.Item(126, 9).Value = 165
$ws.Cells.Item(126, 10).Value = 294954
$ws.Cells.Item(126, 11).Value = 'Hydraulic oil	        1000'
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = 210144
$ws.Cells.Item(126, 14).Value = 'issued'

# Row 127
$ws.Cells.Item(127, 1).Value = '23-03-2018'
$ws.Cells.Item(127, 2).Value = 'v2'
$ws.Cells.Item(127, 3).Value = 'Day/'
$ws.Cells.Item(127, 4).Value = 451
$ws.Cells.Item(127, 5).Value = 95
$ws.Cells.Item(127, 6).Value = 588
$ws.Cells.Item(127, 7).Value = 95
$ws.Cells.Item(127, 8).Value = 95
$ws.Cells.Item(127, 9).Value = 954
$ws.Cells.Item(127, 10).Value = 55860
$ws.Cells.Item(127, 11).Value = 'Hydraulic strainer  250'
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = -34770
$ws.Cells.Item(127, 14).Value = 'issued'

# Row 128
$ws.Cells.Item(128, 1).Value = '13-03-2018'
$ws.Cells.Item(128, 2).Value = 'v3'
$ws.Cells.Item(128, 3).Value = 'Day/'
$ws.Cells.Item(128, 4).Value = 564
$ws.Cells.Item(128, 5).Value = 65
$ws.Cells.Item(128, 6).Value = 64
$ws.Cells.Item(128, 7).Value = 65
$ws.Cells.Item(128, 8).Value = 6548
$ws.Cells.Item(128, 9).Value = 654
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 'Hydraulic oil	        1000'
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 13).Value = -4282392
$ws.Cells.Item(128, 14).Value = 'issued'

# Row 129
$ws.Cells.Item(129, 1).Value = '13-03-2018'
$ws.Cells.Item(129, 2).Value = 'v3'
$ws.Cells.Item(129, 3).Value = 'Day/'
$ws.Cells.Item(129, 4).Value = 65
$ws.Cells.Item(129, 5).Value = 65
$ws.Cells.Item(129, 6).Value = 654
$ws.Cells.Item(129, 7).Value = 65
$ws.Cells.Item(129, 8).Value = 6548
$ws.Cells.Item(129, 9).Value = 654
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 'Hydraulic oil	        1000'
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 13).Value = -4282392
$ws.Cells.Item(129, 14).Value = 'issued'

# Row 130
$c = $ws.Cells.Item(130, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(130, 2).Value = 'v1'
$ws.Cells.Item(130, 3).Value = 'Day/'
$ws.Cells.Item(130, 4).Value = 84
$ws.Cells.Item(130, 5).Value = 14
$ws.Cells.Item(130, 6).Value = 1700
$ws.Cells.Item(130, 7).Value = 14
$ws.Cells.Item(130, 8).Value = 200
$ws.Cells.Item(130, 9).Value = 25
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 13).Value = -5000
$ws.Cells.Item(130, 14).Value = 'issued'

# Row 131
$c = $ws.Cells.Item(131, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(131, 2).Value = 'v1'
$ws.Cells.Item(131, 3).Value = 'Day/Night'
$ws.Cells.Item(131, 4).Value = 14
$ws.Cells.Item(131, 5).Value = 14
$ws.Cells.Item(131, 6).Value = 1730
$ws.Cells.Item(131, 7).Value = 14
$ws.Cells.Item(131, 8).Value = 200
$ws.Cells.Item(131, 9).Value = 25
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 13).Value = -5000
$ws.Cells.Item(131, 14).Value = 'issued'

# Row 132
$c = $ws.Cells.Item(132, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(132, 2).Value = 'v3'
$ws.Cells.Item(132, 3).Value = '/Night'
$ws.Cells.Item(132, 4).Value = 65
$ws.Cells.Item(132, 5).Value = 354
$ws.Cells.Item(132, 6).Value = 2000
$ws.Cells.Item(132, 7).Value = 354
$ws.Cells.Item(132, 8).Value = 84
$ws.Cells.Item(132, 9).Value = 684
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 'Nothing'
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -57456
$ws.Cells.Item(132, 14).Value = 'issued'

# Row 133
$ws.Cells.Item(133, 1).Value = '14-03-2018'
$ws.Cells.Item(133, 2).Value = 'v2'
$ws.Cells.Item(133, 3).Value = 'Day/'
$ws.Cells.Item(133, 4).Value = 95
$ws.Cells.Item(133, 5).Value = 25
$ws.Cells.Item(133, 6).Value = 245
$ws.Cells.Item(133, 7).Value = 25
$ws.Cells.Item(133, 8).Value = 1500
$ws.Cells.Item(133, 9).Value = 245
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 13).Value = -367500
$ws.Cells.Item(133, 14).Value = 'issued'

# Row 134
$ws.Cells.Item(134, 1).Value = '30-03-2018'
$ws.Cells.Item(134, 2).Value = 'v4'
$ws.Cells.Item(134, 3).Value = 'Day/Night'
$ws.Cells.Item(134, 4).Value = 1516
$ws.Cells.Item(134, 5).Value = 2032
$ws.Cells.Item(134, 6).Value = 2245
$ws.Cells.Item(134, 7).Value = 2032
$ws.Cells.Item(134, 8).Value = 66
$ws.Cells.Item(134, 9).Value = 200
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 'Nothing'
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -13200
$ws.Cells.Item(134, 14).Value = 'issued'

# Row 135
$ws.Cells.Item(135, 1).Value = '31-03-2018'
$ws.Cells.Item(135, 2).Value = 'v4'
$ws.Cells.Item(135, 3).Value = 'Day/Night'
$ws.Cells.Item(135, 4).Value = 2032
$ws.Cells.Item(135, 5).Value = 2040
$ws.Cells.Item(135, 6).Value = 2245
$ws.Cells.Item(135, 7).Value = 2040
$ws.Cells.Item(135, 8).Value = 66
$ws.Cells.Item(135, 9).Value = 200
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 'Nothing'
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -13200
$ws.Cells.Item(135, 14).Value = 'issued'

# Row 136
$ws.Cells.Item(136, 1).Value = '31-03-2018'
$ws.Cells.Item(136, 2).Value = 'v4'
$ws.Cells.Item(136, 3).Value = 'Day/Night'
$ws.Cells.Item(136, 4).Value = 2040
$ws.Cells.Item(136, 5).Value = 2040
$ws.Cells.Item(136, 6).Value = 2245
$ws.Cells.Item(136, 7).Value = 2040
$ws.Cells.Item(136, 8).Value = 66
$ws.Cells.Item(136, 9).Value = 200
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 'Hydraulic oil	        1000'
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -13200
$ws.Cells.Item(136, 14).Value = 'issued'

# Row 137
$c = $ws.Cells.Item(137, 1)
$c.NumberFormat = "@"
$c.Value = '07-03-2018'
$c.ClearFormats()
$ws.Cells.Item(137, 2).Value = 'v4'
$ws.Cells.Item(137, 3).Value = '/'
$ws.Cells.Item(137, 4).Value = 2040
$ws.Cells.Item(137, 5).Value = 2048
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 2048
$ws.Cells.Item(137, 8).Value = 2454
$ws.Cells.Item(137, 9).Value = 254
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -623316
$ws.Cells.Item(137, 14).Value = 'issued'

# Row 138
$ws.Cells.Item(138, 1).Value = '24-03-2018'
$ws.Cells.Item(138, 2).Value = 'v3'
$ws.Cells.Item(138, 3).Value = 'Day/Night'
$ws.Cells.Item(138, 4).Value = 354
$ws.Cells.Item(138, 5).Value = 989
$ws.Cells.Item(138, 6).Value = 2024
$ws.Cells.Item(138, 7).Value = 989
$ws.Cells.Item(138, 8).Value = 66
$ws.Cells.Item(138, 9).Value = 200
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 'Engine oil 	        250'
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = -13200
$ws.Cells.Item(138, 14).Value = 'issued'

# Row 139
$ws.Cells.Item(139, 1).Value = '24-03-2018'
$ws.Cells.Item(139, 2).Value = 'v4'
$ws.Cells.Item(139, 3).Value = 'Day/Night'
$ws.Cells.Item(139, 4).Value = 2048
$ws.Cells.Item(139, 5).Value = 1500
$ws.Cells.Item(139, 6).Value = 2024
$ws.Cells.Item(139, 7).Value = 1500
$ws.Cells.Item(139, 8).Value = 66
$ws.Cells.Item(139, 9).Value = 3123
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 'Nothing'
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = -206118
$ws.Cells.Item(139, 14).Value = 'issued'

# Row 140
$ws.Cells.Item(140, 1).Value = '24-03-2018'
$ws.Cells.Item(140, 2).Value = 'v4'
$ws.Cells.Item(140, 3).Value = 'Day/Night'
$ws.Cells.Item(140, 4).Value = 1500
$ws.Cells.Item(140, 5).Value = 1500
$ws.Cells.Item(140, 6).Value = 2024
$ws.Cells.Item(140, 7).Value = 1500
$ws.Cells.Item(140, 8).Value = 66
$ws.Cells.Item(140, 9).Value = 3123
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 'Track motor oil       1000'
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -206118
$ws.Cells.Item(140, 14).Value = 'issued'
